{"js": "// Adds the text \"Peque\u00f1a modificacion\" (underlined) to the empty,\n// underline-formatted paragraph at the end of the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the last paragraph (the empty one that carries the\n// single-underline paragraph mark formatting in the original document).\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert the new run of text into that paragraph.\nconst range = lastParagraph.insertText(\"Peque\u00f1a modificacion\", Word.InsertLocation.end);\nrange.font.underline = Word.UnderlineType.single;\n\nawait context.sync();\n", "ps1": "# Adds the text \"Peque\u00f1a modificacion\" (underlined) to the empty,\n# underline-formatted paragraph at the end of the document body.\n$d = $word.ActiveDocument\n\n# The target is the last paragraph in the document body (the empty\n# paragraph that already carries single-underline formatting).\n$paragraphs = $d.Paragraphs\n$lastParagraph = $paragraphs.Item($paragraphs.Count)\n\n$r = $lastParagraph.Range\n$r.InsertAfter(\"Peque\u00f1a modificacion\")\n$r.Font.Underline = 1\n\nWrite-Host \"Inserted new underlined run into the last paragraph.\"\n"}
